$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new record at row 283, shifting the existing rows 283-338
# (and their formatting) down to 284-339. Excel's row Insert naturally
# copies formats/styles down with the shifted cells, which is exactly
# what the target diff shows (e.g. the date-format style on column D
# following the data down one row at a time, and a new row 339 ending up
# with the values that used to live in row 338).
$ws.Rows.Item(283).Insert()

# Populate the columns that stay identical across every record in this
# block (copied straight from the row directly below, which now holds
# the data that used to be row 283).
$ws.Range("A283").Value2 = $ws.Range("A284").Value2
$ws.Range("B283").Value2 = $ws.Range("B284").Value2
$ws.Range("C283").Value2 = $ws.Range("C284").Value2
$ws.Range("E283").Value2 = $ws.Range("E284").Value2
$ws.Range("F283").Value2 = $ws.Range("F284").Value2
$ws.Range("G283").Value2 = $ws.Range("G284").Value2
$ws.Range("H283").Value2 = $ws.Range("H284").Value2
$ws.Range("I283").Value2 = $ws.Range("I284").Value2
$ws.Range("Q283").Value2 = $ws.Range("Q284").Value2
$ws.Range("R283").Value2 = $ws.Range("R284").Value2

# New record's own values.
$ws.Range("D283").Value2 = 44943
$ws.Range("J283").Value2 = 70
$ws.Range("K283").Value2 = 43000
$ws.Range("L283").Value2 = 45000
$ws.Range("M283").Value2 = 44000
$ws.Range("N283").Value2 = "$/saco 25 kilos"
$ws.Range("O283").Value2 = "Región Metropolitana"
$ws.Range("P283").Value2 = 1760
